# Add data for 2024-11-06 (column BG, rows 2-53).
#
# The sheet has no real conditional formatting - every data cell's
# highlight (white / yellow / light-blue fill, "Meiryo" font) is a plain
# cell style baked in at save time, picked by the value's range:
#   value <  125        -> yellow fill   (existing style used by col D, row 2: ~65535)
#   125 <= value < 140   -> light-blue fill (existing style used by col N, row 2: ~15128749)
#   value >= 140          -> plain/white fill (existing style used by col A: 16777215)
# Column BG currently duplicates column BF; this script overwrites BG2:BG53
# with the real 2024-11-06 values and re-derives the matching style by
# copying the formatting from another cell in the same row that already
# carries the right style (so no new style entries get minted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$COLOR_YELLOW = 65535      # fill used for value < 125
$COLOR_BLUE   = 15128749   # fill used for 125 <= value < 140
$COLOR_WHITE  = 16777215   # fill used for value >= 140 (plain/no highlight)

$xlPasteFormats = -4122

# New 2024/11/06 values for BG2:BG53, in row order.
$newValues = @{
    2  = 139.4;  3  = 158;   4  = 118.8; 5  = 111.4; 6  = 255.5;
    7  = 142.9;  8  = 126.8; 9  = 273;   10 = 148.7; 11 = 131;
    12 = 207.2;  13 = 162.4; 14 = 163.1; 15 = 157;   16 = 296.6;
    17 = 176.6;  18 = 145.5; 19 = 154.2; 20 = 149;   21 = 165.5;
    22 = 128.3;  23 = 121.4; 24 = 160.4; 25 = 146.1; 26 = 139.2;
    27 = 260.7;  28 = 228.3; 29 = 141.2; 30 = 145.3; 31 = 173.9;
    32 = 195;    33 = 155.2; 34 = 240.8; 35 = 142.6; 36 = 127.4;
    37 = 163.7;  38 = 143.8; 39 = 130.2; 40 = 243.6; 41 = 197.9;
    42 = 160.1;  43 = 153.6; 44 = 161;   45 = 144.5; 46 = 159;
    47 = 238.9;  48 = 151.3; 49 = 223;   50 = 153.5; 51 = 161.1;
    52 = 199.2;  53 = 164.9
}

function Get-TargetColor($value) {
    if ($value -lt 125) { return $COLOR_YELLOW }
    elseif ($value -lt 140) { return $COLOR_BLUE }
    else { return $COLOR_WHITE }
}

function Find-TemplateCell($ws, $row, $wantColor) {
    # Look across the row (columns B..BF; BG is the column we're about to
    # overwrite) for an existing cell whose fill already matches the style
    # we need, so copying its format reuses the workbook's existing style
    # instead of minting a new one.
    for ($c = 2; $c -le 58; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        if ($cell.Interior.Color -eq $wantColor) {
            return $cell
        }
    }
    # Column A (the machine-number column) is always the plain/white style.
    return $ws.Cells.Item($row, 1)
}

foreach ($row in 2..53) {
    $value = $newValues[$row]
    $wantColor = Get-TargetColor $value

    $template = Find-TemplateCell $ws $row $wantColor
    $target = $ws.Cells.Item($row, 59)   # column BG

    $template.Copy() | Out-Null
    $target.PasteSpecial($xlPasteFormats) | Out-Null
    $target.Value = $value
}

$excel.CutCopyMode = $false
Write-Output "Updated BG2:BG53 with 2024-11-06 data"
